$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Version: bump patch build date from 202312 to 202406
$ws1.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Title: falls back to Name value (CDACompressionAlgorithm) instead of "CompressionAlgorithm"
$ws1.Range("B5").Value = "CDACompressionAlgorithm"

# Date: refreshed publication date/time
$ws1.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact: proper HL7 Structured Documents contact info instead of placeholder
$ws1.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"
